$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "23.460.54"
$ws.Cells.Item(2, 5).Value = "  +1.17%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.639.15"
$ws.Cells.Item(3, 5).Value = "  +2.44%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 5).Value = "  +0.05%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "306.64"
$ws.Cells.Item(6, 5).Value = "  +1.14%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.3763"
$ws.Cells.Item(7, 5).Value = "  -0.44%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "52.28"
$ws.Cells.Item(8, 5).Value = "  +0.21%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.3644"
$ws.Cells.Item(9, 5).Value = "  +0.89%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.264"
$ws.Cells.Item(10, 5).Value = "  -0.04%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.08154"
$ws.Cells.Item(11, 5).Value = "  +0.47%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.001"
$ws.Cells.Item(12, 5).Value = "  -0.02%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "22.92"
$ws.Cells.Item(13, 5).Value = "  +1.14%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.634"
$ws.Cells.Item(14, 5).Value = "  +0.90%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.00001275"
$ws.Cells.Item(15, 5).Value = "  +2.43%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "7.368"
$ws.Cells.Item(16, 5).Value = "  -0.41%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "1.640.55"
$ws.Cells.Item(17, 5).Value = "  +2.48%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "94.59"
$ws.Cells.Item(18, 5).Value = "  +0.64%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06927"
$ws.Cells.Item(19, 5).Value = "  +0.57%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "18.15"
$ws.Cells.Item(20, 5).Value = "  +0.51%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.550"
$ws.Cells.Item(21, 5).Value = "  +0.13%  "
$ws.Cells.Item(22, 5).Value = "  +0.02%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "23.457.97"
$ws.Cells.Item(23, 5).Value = "  +1.20%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "12.79"
$ws.Cells.Item(24, 5).Value = "  -1.26%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "3.101"
$ws.Cells.Item(25, 5).Value = "  +3.38%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.421"
$ws.Cells.Item(26, 5).Value = "  +1.59%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "21.24"
$ws.Cells.Item(27, 5).Value = "  +0.01%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "150.74"
$ws.Cells.Item(28, 5).Value = "  +0.86%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "5.362"
$ws.Cells.Item(29, 5).Value = "  +2.24%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "135.28"
$ws.Cells.Item(30, 5).Value = "  +1.11%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "2.302"
$ws.Cells.Item(31, 5).Value = "  -3.43%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.824.78"
$ws.Cells.Item(32, 5).Value = "  +2.53%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "6.794"
$ws.Cells.Item(33, 5).Value = "  +0.01%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.9662"
$ws.Cells.Item(34, 5).Value = "  +0.08%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.02822"
$ws.Cells.Item(35, 5).Value = "  +4.32%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "10.32"
$ws.Cells.Item(36, 5).Value = "  +0.56%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.07301"
$ws.Cells.Item(37, 5).Value = "  -2.51%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.2525"
$ws.Cells.Item(38, 5).Value = "  +0.89%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.08837"
$ws.Cells.Item(39, 5).Value = "  +0.39%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "6.121"
$ws.Cells.Item(40, 5).Value = "  +0.67%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.377"
$ws.Cells.Item(41, 5).Value = "  +1.14%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.7106"
$ws.Cells.Item(42, 5).Value = "  +0.13%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "12.52"
$ws.Cells.Item(43, 5).Value = "  +0.44%  "
$ws.Cells.Item(44, 5).Value = "  +4.42%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.6543"
$ws.Cells.Item(45, 5).Value = "  +0.36%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.338"
$ws.Cells.Item(46, 5).Value = "  +1.20%  "
$ws.Cells.Item(47, 5).Value = "  +0.17%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "4.022"
$ws.Cells.Item(48, 5).Value = "  +0.28%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.07959"
$ws.Cells.Item(49, 5).Value = "  +0.01%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "128.85"
$ws.Cells.Item(50, 5).Value = "  -2.41%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.201"
$ws.Cells.Item(51, 5).Value = "  +0.06%  "
